$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving the default (unstyled) cell
# format -- matters for numeric-looking strings (e.g. "576.82") which Excel
# would otherwise silently coerce into a real number on assignment.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "61.648.57"
Set-TextValue $ws.Range("E2") "  -2.01%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.395.55"
Set-TextValue $ws.Range("E3") "  -1.51%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.21%  "

# Row 5
Set-TextValue $ws.Range("D5") "576.82"
Set-TextValue $ws.Range("E5") "  +0.21%  "

# Row 6
Set-TextValue $ws.Range("D6") "152.18"
Set-TextValue $ws.Range("E6") "  +3.05%  "

# Row 7
Set-TextValue $ws.Range("E7") "  +0.13%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.484"
Set-TextValue $ws.Range("E8") "  +1.46%  "

# Row 9
Set-TextValue $ws.Range("D9") "8.04"
Set-TextValue $ws.Range("E9") "  +3.45%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.124"
Set-TextValue $ws.Range("E10") "  +0.35%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.420"
Set-TextValue $ws.Range("E11") "  +3.49%  "

# Row 12
Set-TextValue $ws.Range("D12") "3.991.68"
Set-TextValue $ws.Range("E12") "  -1.17%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.130"
Set-TextValue $ws.Range("E13") "  +1.26%  "

# Row 14
Set-TextValue $ws.Range("D14") "28.56"
Set-TextValue $ws.Range("E14") "  -0.42%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D15") "3.428.00"
Set-TextValue $ws.Range("E15") "  -0.66%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D16") "0.0000172"
Set-TextValue $ws.Range("E16") "  +0.63%  "

# Row 17
Set-TextValue $ws.Range("D17") "61.871.30"
Set-TextValue $ws.Range("E17") "  -1.72%  "

# Row 18
Set-TextValue $ws.Range("D18") "6.42"
Set-TextValue $ws.Range("E18") "  +0.59%  "

# Row 19
Set-TextValue $ws.Range("D19") "14.35"
Set-TextValue $ws.Range("E19") "  +0.21%  "

# Row 20
Set-TextValue $ws.Range("D20") "9.01"
Set-TextValue $ws.Range("E20") "  -1.12%  "

# Row 21
Set-TextValue $ws.Range("D21") "377.14"
Set-TextValue $ws.Range("E21") "  -1.89%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.568"
Set-TextValue $ws.Range("E22") "  +1.91%  "

# Row 23
Set-TextValue $ws.Range("D23") "75.89"
Set-TextValue $ws.Range("E23") "  +1.85%  "

# Row 24
Set-TextValue $ws.Range("E24") "  +0.02%  "

# Row 25
Set-TextValue $ws.Range("D25") "3.561.12"
Set-TextValue $ws.Range("E25") "  -0.49%  "

# Row 26
Set-TextValue $ws.Range("D26") "0.0000110"
Set-TextValue $ws.Range("E26") "  -3.29%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.177"
Set-TextValue $ws.Range("E27") "  -2.56%  "

# Row 28
Set-TextValue $ws.Range("D28") "7.57"
Set-TextValue $ws.Range("E28") "  -0.02%  "

# Row 29
Set-TextValue $ws.Range("D29") "1.00"
Set-TextValue $ws.Range("E29") "  -0.02%  "

# Row 30
Set-TextValue $ws.Range("D30") "2.12"
Set-TextValue $ws.Range("E30") "  +1.57%  "

# Row 31
Set-TextValue $ws.Range("D31") "7.77"
Set-TextValue $ws.Range("E31") "  -2.53%  "

# Row 32
Set-TextValue $ws.Range("E32") "  +0.10%  "

# Row 33
Set-TextValue $ws.Range("D33") "23.05"
Set-TextValue $ws.Range("E33") "  -0.89%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.30"
Set-TextValue $ws.Range("E34") "  +2.00%  "

# Row 35
Set-TextValue $ws.Range("D35") "5.43"
Set-TextValue $ws.Range("E35") "  +2.23%  "

# Row 36
Set-TextValue $ws.Range("D36") "1.58"
Set-TextValue $ws.Range("E36") "  -1.48%  "

# Row 37
Set-TextValue $ws.Range("D37") "6.89"
Set-TextValue $ws.Range("E37") "  -2.09%  "

# Row 38
Set-TextValue $ws.Range("D38") "168.82"
Set-TextValue $ws.Range("E38") "  -0.26%  "

# Row 39
Set-TextValue $ws.Range("D39") "30.53"
Set-TextValue $ws.Range("E39") "  -4.21%  "

# Row 40
$ws.Range("B40").Value = "RenzoRestakedETH"
$ws.Range("C40").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue $ws.Range("D40") "3.450.87"
Set-TextValue $ws.Range("E40") "  -0.93%  "

# Row 41
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D41") "0.0773"
Set-TextValue $ws.Range("E41") "  +1.12%  "

# Row 42
Set-TextValue $ws.Range("D42") "42.53"
Set-TextValue $ws.Range("E42") "  +0.34%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.773"
Set-TextValue $ws.Range("E43") "  -2.32%  "

# Row 44
Set-TextValue $ws.Range("D44") "4.38"
Set-TextValue $ws.Range("E44") "  +1.24%  "

# Row 45
Set-TextValue $ws.Range("D45") "1.65"
Set-TextValue $ws.Range("E45") "  -3.87%  "

# Row 46
Set-TextValue $ws.Range("D46") "1.14"
Set-TextValue $ws.Range("E46") "  -3.66%  "

# Row 47
Set-TextValue $ws.Range("D47") "2.533.91"
Set-TextValue $ws.Range("E47") "  -1.18%  "

# Row 48
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D48") "6.81"
Set-TextValue $ws.Range("E48") "  -0.95%  "

# Row 49
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D49") "23.00"
Set-TextValue $ws.Range("E49") "  +2.84%  "

# Row 50
Set-TextValue $ws.Range("D50") "1.01"
Set-TextValue $ws.Range("E50") "  +0.66%  "

# Row 51
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D51") "0.0263"
Set-TextValue $ws.Range("E51") "  -1.43%  "
